$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.618.11'
$ws.Range('E2').Value = '  +0.20%  '

$ws.Range('D3').Value = '1.812.40'
$ws.Range('E3').Value = '  -0.02%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.40'
$ws.Range('E5').Value = '  -1.00%  '

$ws.Range('E6').Value = '  +3.46%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '37.80'
$ws.Range('E8').Value = '  +8.08%  '

$ws.Range('E9').Value = '  -3.53%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0680'
$ws.Range('E10').Value = '  -2.52%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('E11').Value = '  +1.58%  '

$ws.Range('D12').Value = '2.074.57'
$ws.Range('E12').Value = '  +0.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.31'
$ws.Range('E13').Value = '  +0.56%  '

$ws.Range('D14').Value = '1.829.73'
$ws.Range('E14').Value = '  +0.84%  '

$ws.Range('E15').Value = '  -2.45%  '

$ws.Range('D16').Value = '34.586.84'
$ws.Range('E16').Value = '  +0.24%  '

$ws.Range('E17').Value = '  -1.91%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.77'
$ws.Range('E18').Value = '  -0.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.53'
$ws.Range('E19').Value = '  -0.79%  '

$ws.Range('E20').Value = '  -2.76%  '

$ws.Range('E21').Value = '  -1.87%  '

$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  -0.78%  '

$ws.Range('E24').Value = '  +4.68%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.06'
$ws.Range('E25').Value = '  -0.18%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.94'
$ws.Range('E26').Value = '  -1.87%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.35'
$ws.Range('E27').Value = '  +2.95%  '

$ws.Range('E28').Value = '  +1.50%  '

$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.93'
$ws.Range('E30').Value = '  -2.80%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.81'
$ws.Range('E31').Value = '  -1.52%  '

$ws.Range('E32').Value = '  -1.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0524'
$ws.Range('E33').Value = '  -2.35%  '

$ws.Range('E34').Value = '  -1.05%  '

$ws.Range('D35').Value = '1.365.95'
$ws.Range('E35').Value = '  -2.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.656'
$ws.Range('E36').Value = '  -3.66%  '

$ws.Range('E37').Value = '  -0.64%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0188'
$ws.Range('E38').Value = '  -1.67%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.34'
$ws.Range('E39').Value = '  -5.35%  '

$ws.Range('E40').Value = '  +8.00%  '

$ws.Range('E41').Value = '  +1.38%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '81.14'

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.78'
$ws.Range('E43').Value = '  -1.70%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.941'
$ws.Range('E44').Value = '  -2.52%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.15'
$ws.Range('E45').Value = '  +6.25%  '

$ws.Range('E46').Value = '  -2.30%  '

$ws.Range('D47').Value = '1.974.93'
$ws.Range('E47').Value = '  +0.06%  '

$ws.Range('E48').Value = '  -3.17%  '

$ws.Range('E49').Value = '  +0.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.78'
$ws.Range('E50').Value = '  -2.49%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0617'
$ws.Range('E51').Value = '  +1.60%  '
